# Risolto bug HMR e aggiunti altri log per il debug
#
# The "Fallimenti per Fragilità" (fragility-failure) count for the
# LLM / robula row (row 6) was wrong: it should be 3, not 6. Updating it
# causes the dependent formulas (D6, G6, B20, D20) to recalculate
# automatically. Also move the active selection to F6, the cell that was
# edited, matching what a user would do after fixing the value by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 3

$ws.Range("F6").Select()
